{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\n// The document has 4 paragraphs:\n//   1. \"Basic if demonstration :\"\n//   2. \"The ELSE paragraph.\"\n//   3. \"End of demonstration.\"\n//   4. \"\" (empty trailing paragraph)\n//\n// Target (\"Improved Block error reporting\"): drop the stray ELSE /\n// end-of-demo paragraphs and the empty trailing paragraph, replacing them\n// with a single bold red error message appended to the end of the first\n// paragraph.\n\nconst paragraphs = body.paragraphs.items;\nconst firstParagraph = paragraphs[0];\n\nconst elseParagraph = paragraphs.find((p) => p.text.trim() === \"The ELSE paragraph.\");\nconst endParagraph = paragraphs.find((p) => p.text.trim() === \"End of demonstration.\");\nconst lastParagraph = paragraphs[paragraphs.length - 1];\n\n// Delete the \"ELSE paragraph\" and \"End of demonstration\" paragraphs\n// entirely.\nelseParagraph.delete();\nendParagraph.delete();\nawait context.sync();\n\n// Append the error message to the end of the first paragraph, formatted\n// bold and red.\nconst errorRange = firstParagraph.insertText(\n  \"Invalid if statement: Unexpected tag EOF missing [ENDIF]\",\n  Word.InsertLocation.end\n);\nerrorRange.font.bold = true;\nerrorRange.font.color = \"#FF0000\";\nawait context.sync();\n\n// Remove the now-orphaned empty trailing paragraph (former paragraph 4).\nlastParagraph.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document has 4 paragraphs:\n#   1. \"Basic if demonstration :\"\n#   2. \"The ELSE paragraph.\"\n#   3. \"End of demonstration.\"\n#   4. \"\" (empty trailing paragraph)\n#\n# Target (\"Improved Block error reporting\"): drop the stray ELSE /\n# end-of-demo paragraphs and the empty trailing paragraph, replacing them\n# with a single bold red error message appended to the end of the first\n# paragraph.\n\n$elseParaIndex = -1\n$endParaIndex = -1\n$lastParaIndex = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"The ELSE paragraph.\") { $elseParaIndex = $i }\n    if ($t -eq \"End of demonstration.\") { $endParaIndex = $i }\n}\n\n$firstParagraph = $d.Paragraphs.Item(1)\n\n# Delete the \"ELSE paragraph\" and \"End of demonstration\" paragraphs\n# entirely (they are adjacent, so remove them as one contiguous range).\n$delRange = $d.Range($d.Paragraphs.Item($elseParaIndex).Range.Start, $d.Paragraphs.Item($endParaIndex).Range.End)\n$delRange.Delete()\n\n# Append the error message to the end of the first paragraph, just before\n# its paragraph mark, and format it bold + red.\n$insertPoint = $d.Range($firstParagraph.Range.End - 1, $firstParagraph.Range.End - 1)\n$insertPoint.InsertAfter(\"Invalid if statement: Unexpected tag EOF missing [ENDIF]\")\n$insertPoint.Font.Bold = $true\n$insertPoint.Font.Color = 255\n\n# Remove the now-orphaned empty trailing paragraph (former paragraph 4).\n$d.Paragraphs.Item(2).Range.Delete()\n"}
